# Season-record columns: Wins / Losses / Ties appended after the existing
# team-stats columns (A:AC), mirroring the header style used by the rest
# of row 1 and filling every player row (2-47) with the team's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
# Copy the formatting (bold font, border, center/top alignment) from an
# existing header cell onto the three new header cells, then set text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-47): the team's season record, same for every player ---
$ws.Range("AD2:AD47").Value = 96
$ws.Range("AE2:AE47").Value = 66
$ws.Range("AF2:AF47").Value = 0
